$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E2").Value = 1194.9707055
$wsSchedule.Range("F2").Value = 19.75811351686508

# --- Sheet: Detailed ---
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsDetailed.Range("B16").Value = 26.82191

$wsDetailed.Range("B17").Value = 22.12431
$wsDetailed.Range("C17").Value = "historical"

$wsDetailed.Range("B18").Value = 20.94111
$wsDetailed.Range("C18").Value = "historical"

$wsDetailed.Range("B19").Value = 32.11392
$wsDetailed.Range("B20").Value = 32.61874
$wsDetailed.Range("B21").Value = 22.07
$wsDetailed.Range("B22").Value = 22.07
$wsDetailed.Range("B23").Value = 5.17355
$wsDetailed.Range("B24").Value = 5.26655
$wsDetailed.Range("B25").Value = 22.78583
$wsDetailed.Range("B26").Value = 22.07
$wsDetailed.Range("B27").Value = 25.8407
$wsDetailed.Range("B28").Value = 22.07
$wsDetailed.Range("B29").Value = 22.07
$wsDetailed.Range("B30").Value = 26.61817
$wsDetailed.Range("B33").Value = 22.07
$wsDetailed.Range("B34").Value = 22.07
$wsDetailed.Range("B35").Value = 0.51
$wsDetailed.Range("B36").Value = 12.44773
$wsDetailed.Range("B37").Value = 30.08846
$wsDetailed.Range("B38").Value = 46.25725
$wsDetailed.Range("B39").Value = 57.00963
$wsDetailed.Range("B41").Value = 58.29804
$wsDetailed.Range("B42").Value = 57.09
$wsDetailed.Range("B44").Value = 57.09
